$d = $word.ActiveDocument
$d.Content.Find.Execute("2026-02-01 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-02-02 Monday", 2)
$d.Content.Find.Execute("40÷2=20, 0", $true, $false, $false, $false, $false, $true, 1, $false, "62÷3=20, 2", 2)
$d.Content.Find.Execute("34÷6=5, 4", $true, $false, $false, $false, $false, $true, 1, $false, "63÷6=10, 3", 2)
$d.Content.Find.Execute("47÷3=15, 2", $true, $false, $false, $false, $false, $true, 1, $false, "65÷7=9, 2", 2)
$d.Content.Find.Execute("53÷3=17, 2", $true, $false, $false, $false, $false, $true, 1, $false, "52÷4=13, 0", 2)
$d.Content.Find.Execute("84÷4=21, 0", $true, $false, $false, $false, $false, $true, 1, $false, "24÷9=2, 6", 2)
$d.Content.Find.Execute("32÷4=8, 0", $true, $false, $false, $false, $false, $true, 1, $false, "58÷2=29, 0", 2)
$d.Content.Find.Execute("85÷3=28, 1", $true, $false, $false, $false, $false, $true, 1, $false, "53÷6=8, 5", 2)
$d.Content.Find.Execute("65÷4=16, 1", $true, $false, $false, $false, $false, $true, 1, $false, "14÷3=4, 2", 2)
$d.Content.Find.Execute("42÷3=14, 0", $true, $false, $false, $false, $false, $true, 1, $false, "87÷3=29, 0", 2)
$d.Content.Find.Execute("46÷9=5, 1", $true, $false, $false, $false, $false, $true, 1, $false, "54÷9=6, 0", 2)
$d.Content.Find.Execute("71÷9=7, 8", $true, $false, $false, $false, $false, $true, 1, $false, "26÷7=3, 5", 2)
$d.Content.Find.Execute("45÷8=5, 5", $true, $false, $false, $false, $false, $true, 1, $false, "42÷9=4, 6", 2)
$d.Content.Find.Execute("94÷6=15, 4", $true, $false, $false, $false, $false, $true, 1, $false, "80÷2=40, 0", 2)
$d.Content.Find.Execute("11÷9=1, 2", $true, $false, $false, $false, $false, $true, 1, $false, "91÷5=18, 1", 2)
$d.Content.Find.Execute("20÷5=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "91÷9=10, 1", 2)
$d.Content.Find.Execute("69÷8=8, 5", $true, $false, $false, $false, $false, $true, 1, $false, "82÷7=11, 5", 2)
$d.Content.Find.Execute("98÷2=49, 0", $true, $false, $false, $false, $false, $true, 1, $false, "82÷6=13, 4", 2)
$d.Content.Find.Execute("68÷7=9, 5", $true, $false, $false, $false, $false, $true, 1, $false, "45÷9=5, 0", 2)
$d.Content.Find.Execute("60÷4=15, 0", $true, $false, $false, $false, $false, $true, 1, $false, "68÷4=17, 0", 2)
$d.Content.Find.Execute("17÷7=2, 3", $true, $false, $false, $false, $false, $true, 1, $false, "63÷3=21, 0", 2)
$d.Content.Find.Execute("70÷2=35, 0", $true, $false, $false, $false, $false, $true, 1, $false, "57÷8=7, 1", 2)
$d.Content.Find.Execute("49÷8=6, 1", $true, $false, $false, $false, $false, $true, 1, $false, "66÷6=11, 0", 2)
$d.Content.Find.Execute("27÷2=13, 1", $true, $false, $false, $false, $false, $true, 1, $false, "84÷2=42, 0", 2)
$d.Content.Find.Execute("60÷8=7, 4", $true, $false, $false, $false, $false, $true, 1, $false, "24÷9=2, 6", 2)
$d.Content.Find.Execute("58÷6=9, 4", $true, $false, $false, $false, $false, $true, 1, $false, "44÷8=5, 4", 2)
Write-Host "done"
